$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dates (column A) and tide values (column B) for rows 2-26
$dates = @(
    "20000106","20000120","20000205","20000219","20000305",
    "20000319","20000404","20000418","20000503","20000517",
    "20000602","20000616","20000701","20000716","20000730",
    "20000814","20000829","20000913","20000927","20001013",
    "20001027","20001111","20001125","20001211","20001225"
)

$values = @(
    4.905,6.66,4.836,6.05,5.276,
    6.178,5.522,6.657,5.594,5.538,
    6.325,5.476,6.306,5.666,6.424,
    6.037,6.447,5.446,6.841,5.331,
    6.066,6.526,6.759,6.375,4.747
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    # Force the cell to be treated as text so the numeric-looking date
    # string is not auto-converted into a number, then restore the
    # original (default) cell style so no extra formatting is left behind.
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the summary rows: "Median:" value, relabel "Average:" to "mean:",
# and update the mean value.
$ws.Range("B28").Value = 6.0435
$ws.Range("A29").Value = "mean:"
$ws.Range("B29").Value = 5.91952
